$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Price column values are forced to text (matching the original inline string
# cell type) since many look like numbers (e.g. "1.001", "0.9999") and would
# otherwise be auto-converted by Excel; the style is reset back to Normal so
# no extra formatting is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.215.25'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.58%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.859.09'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.94%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6994'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.50%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07824'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.62%  '

$ws.Range("E9").Value = '  -0.94%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.02'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.76%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07805'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.99%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.858.74'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.42%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.131'
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '92.10'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6932'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.95%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.594'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008540'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.06%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.234.15'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.56%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '248.28'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.84%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.112.56'
$ws.Range("D20").Style = "Normal"

$ws.Range("E21").Value = '  -3.52%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.580'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.12%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.01%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1538'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.80'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.66%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.923'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.58%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.58'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.76%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.575'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.55%  '

$ws.Range("E30").Value = '  -3.12%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.244'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.21%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.209'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.08%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05245'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.51%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7604'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.875'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.62%  '

$ws.Range("E36").Value = '  -0.01%  '

$ws.Range("E37").Value = '  -0.17%  '

$ws.Range("E38").Value = '  -1.88%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.236.54'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.43%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.742'
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9019'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.59%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '109.99'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.71%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.885'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.89%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9998'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.09%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '68.35'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.92%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.011.13'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.13%  '

$ws.Range("E47").Value = '  -4.10%  '

$ws.Range("E48").Value = '  -0.34%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.529'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.03%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.767'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.41%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4257'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.07%  '
